# Apply scraping update for 141 horarios workbook: 31/12 15:41:06 run
# adding newly-scraped rows to LP1912 and 6203-6173 sheets, and refreshing
# the "Última actualización" timestamp on all three sheets.

$wb = $excel.ActiveWorkbook

$newTimestamp = "Última actualización: 31/12/2025 15:41:06"

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = $newTimestamp
$ws1.Range("A3").Value = "Total filas: 1053"

$data1 = @(
    @(1041, "15:40:55", "15:57", "16_SANTA ANA", 17, "LP1912", "31/12/2025"),
    @(1042, "15:40:55", "16:01", "15_ABASTO", 21, "LP1912", "31/12/2025"),
    @(1043, "15:40:55", "16:03", "23_HERNANDEZ", 23, "LP1912", "31/12/2025"),
    @(1044, "15:40:55", "16:09", "16_SANTA ANA", 29, "LP1912", "31/12/2025"),
    @(1045, "15:40:55", "16:13", "10_OLMOS", 33, "LP1912", "31/12/2025"),
    @(1046, "15:40:55", "16:21", "16_SANTA ANA", 41, "LP1912", "31/12/2025"),
    @(1047, "15:40:55", "16:24", "11_ETCHEVERRY", 44, "LP1912", "31/12/2025"),
    @(1048, "15:40:55", "16:31", "16_P MOR-SANTA ANA", 51, "LP1912", "31/12/2025"),
    @(1049, "15:40:55", "16:31", "23_HERNANDEZ", 51, "LP1912", "31/12/2025"),
    @(1050, "15:40:55", "16:36", "17X38_ROMERO", 56, "LP1912", "31/12/2025"),
    @(1051, "15:40:55", "16:53", "10_OLMOS", 73, "LP1912", "31/12/2025"),
    @(1052, "15:40:55", "17:04", "14_ABASTO", 84, "LP1912", "31/12/2025"),
    @(1053, "15:40:55", "17:06", "10_OLMOS", 86, "LP1912", "31/12/2025"),
    @(1054, "15:40:55", "17:07", "15_ABASTO", 87, "LP1912", "31/12/2025")
)

foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = ""
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215 (only the refreshed timestamp changes)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $newTimestamp

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = $newTimestamp
$ws3.Range("A3").Value = "Total filas: 130"

$data3 = @(
    @(130, "31/12/2025", "15:41:01", "15:46", "215C_LA PLATA", 5, "L6203"),
    @(131, "31/12/2025", "15:41:01", "16:59", "215C_LA PLATA", 78, "L6203")
)

foreach ($row in $data3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = ""
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}
